$d = $word.ActiveDocument

# Remove the standalone "Me " run that used to precede the lawyerName
# placeholder, e.g. "... avec Me {lawyerName}." -> "... avec {lawyerName}."
$d.Content.Find.Execute("Me ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)
